# Daily data refresh for the "Pais" (COVID-19 by-country) sheet.
# - Updates the "Datos actualizados ..." timestamp in A1.
# - Writes refreshed Casos totales/Nuevos casos/Casos activos/Recuperados/
#   Casos criticos/Muertes hoy/Muertes figures (columns B-H) for every
#   country row whose numbers moved since the last refresh.
# - Turquia/Indonesia (rows 23-24) and Chequia/Republica Dominicana
#   (rows 39-40) swap ranking order as part of this refresh, so their
#   country names (column A) are rewritten along with their numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 12 de Octubre de 2020 a las 18:08"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 8000852
$ws.Range("C4").Value = 8854
$ws.Range("D4").Value = 5138374
$ws.Range("E4").Value = 2642681
$ws.Range("G4").Value = 102
$ws.Range("H4").Value = 219797

# Row 5: India
$ws.Range("B5").Value = 7160805
$ws.Range("C5").Value = 41505
$ws.Range("D5").Value = 6203130
$ws.Range("E5").Value = 848008
$ws.Range("G5").Value = 483
$ws.Range("H5").Value = 109667

# Row 6: Brasil
$ws.Range("B6").Value = 5096209
$ws.Range("C6").Value = 1230
$ws.Range("E6").Value = 475489
$ws.Range("G6").Value = 49
$ws.Range("H6").Value = 150555

# Row 15: Reino Unido
$ws.Range("B15").Value = 617688
$ws.Range("C15").Value = 13972
$ws.Range("G15").Value = 50
$ws.Range("H15").Value = 42875

# Row 17: Chile
$ws.Range("B17").Value = 482888
$ws.Range("C17").Value = 1517
$ws.Range("D17").Value = 454484
$ws.Range("E17").Value = 15025
$ws.Range("G17").Value = 61
$ws.Range("H17").Value = 13379

# Row 20: Italia
$ws.Range("B20").Value = 359569
$ws.Range("C20").Value = 4619
$ws.Range("D20").Value = 240600
$ws.Range("E20").Value = 82764
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = 36205

# Row 23: Turquia
$ws.Range("A23").Value = "Turquia"
$ws.Range("B23").Value = 337147
$ws.Range("C23").Value = 1614
$ws.Range("D23").Value = 295658
$ws.Range("E23").Value = 32594
$ws.Range("G23").Value = 58
$ws.Range("H23").Value = 8895

# Row 24: Indonesia
$ws.Range("A24").Value = "Indonesia"
$ws.Range("B24").Value = 336716
$ws.Range("C24").Value = 3267
$ws.Range("D24").Value = 258519
$ws.Range("E24").Value = 66262
$ws.Range("G24").Value = 91
$ws.Range("H24").Value = 11935

# Row 25: Alemania
$ws.Range("B25").Value = 328736
$ws.Range("C25").Value = 2445
$ws.Range("E25").Value = 42128
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 9708

# Row 29: Canada
$ws.Range("B29").Value = 182707
$ws.Range("C29").Value = 843
$ws.Range("D29").Value = 154237
$ws.Range("E29").Value = 18845
$ws.Range("G29").Value = 12
$ws.Range("H29").Value = 9625

# Row 39: Chequia
$ws.Range("A39").Value = "Chequia"
$ws.Range("B39").Value = 119007
$ws.Range("C39").Value = 1897
$ws.Range("D39").Value = 56440
$ws.Range("E39").Value = 61522
$ws.Range("G39").Value = 58
$ws.Range("H39").Value = 1045

# Row 40: Republica Dominicana
$ws.Range("A40").Value = "Republica Dominicana"
$ws.Range("B40").Value = 118843
$ws.Range("C40").Value = 366
$ws.Range("D40").Value = 94532
$ws.Range("E40").Value = 22132
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 2179

# Row 48: Guatemala
$ws.Range("B48").Value = 97826
$ws.Range("C48").Value = 111
$ws.Range("D48").Value = 87016
$ws.Range("E48").Value = 7423
$ws.Range("G48").Value = 3
$ws.Range("H48").Value = 3387

# Row 58: Suiza
$ws.Range("D58").Value = 49500
$ws.Range("E58").Value = 12840
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 2096

# Row 59: Moldavia
$ws.Range("B59").Value = 62618
$ws.Range("C59").Value = 467
$ws.Range("D59").Value = 44728
$ws.Range("E59").Value = 16412
$ws.Range("G59").Value = 17
$ws.Range("H59").Value = 1478

# Row 62: Singapur
$ws.Range("D62").Value = 57728
$ws.Range("E62").Value = 125

# Row 84: Jordania
$ws.Range("B84").Value = 26073
$ws.Range("C84").Value = 1147
$ws.Range("D84").Value = 6219
$ws.Range("E84").Value = 19647
$ws.Range("G84").Value = 16
$ws.Range("H84").Value = 207

# Row 87: Grecia
$ws.Range("B87").Value = 22652
$ws.Range("C87").Value = 294
$ws.Range("E87").Value = 12207
$ws.Range("G87").Value = 7
$ws.Range("H87").Value = 456

# Row 89: Republica de Macedonia
$ws.Range("B89").Value = 21113
$ws.Range("C89").Value = 176
$ws.Range("D89").Value = 16301
$ws.Range("E89").Value = 4015
$ws.Range("G89").Value = 5
$ws.Range("H89").Value = 797

# Row 99: Montenegro
$ws.Range("B99").Value = 14050
$ws.Range("C99").Value = 181
$ws.Range("D99").Value = 9920
$ws.Range("E99").Value = 3919
$ws.Range("G99").Value = 9
$ws.Range("H99").Value = 211

# Row 106: Consejo Danes para los Refugiados
$ws.Range("B106").Value = 10868
$ws.Range("C106").Value = 17
$ws.Range("D106").Value = 10249
$ws.Range("E106").Value = 343

# Row 111: Luxemburgo
$ws.Range("B111").Value = 9731
$ws.Range("C111").Value = 9
$ws.Range("D111").Value = 8234
$ws.Range("E111").Value = 1364
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 133

# Row 116: Jamaica
$ws.Range("B116").Value = 7813
$ws.Range("C116").Value = 95
$ws.Range("D116").Value = 3237
$ws.Range("E116").Value = 4430
$ws.Range("G116").Value = 7
$ws.Range("H116").Value = 146

# Row 122: Cuba
$ws.Range("B122").Value = 6000
$ws.Range("C122").Value = 22
$ws.Range("D122").Value = 5574
$ws.Range("E122").Value = 303

# Row 177: Burundi
$ws.Range("B177").Value = 525
$ws.Range("C177").Value = 1
$ws.Range("E177").Value = 52
